$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados", $true, $false, $false, $false, $false, $true, 1, $false, "1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. ^l2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. ^l3. Conceitos de Interface^l4. Compósitos de matriz metálica: características e processos de fabricação. ^l5. Compósitos de matriz cerâmica: características e processos de fabricação. ^l6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. ^l7. Compósitos nanoestruturados. ^l8. Compósitos Naturais. ^l9. Compósitos Híbridos ^l10. Mecânica de estruturas reforçadas. ^lConteúdo prático: ^l1. Caracterização e análise de compósitos de matriz metálica. ^l2. Preparação e caracterização de compósitos de matriz polimérica.^l(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) ^l3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados", 2)
Write-Output "Found: $found"
